$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new log entry row (row 52)
$ws.Range("A52").NumberFormat = $ws.Range("A51").NumberFormat
$ws.Range("A52").Value = "2/5/2025"
$ws.Range("B52").Value = "library storyline"
$ws.Range("C52").Value = 4

# Update the active selection to match the new state
$ws.Range("C46").Select()
